$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the player name
$ws.Name = "Chris Jordan"

# Force the whole used range to be stored as text (all values in the
# source data are text, even the numeric-looking ones such as runs,
# balls, etc. - matches the ignoredError numberStoredAsText flag)
$ws.Range("A1:M3").NumberFormat = "@"

# Header row (a new "matchNo" column was inserted before "teamName")
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 - new match inserted before the previously existing row
$row2 = @("21st","Punjab Kings","Chris Jordan","b Prasidh Krishna","30","18","1","3","166.66","Kolkata Knight Riders","Ahmedabad","April 26","KKR won by 5 wickets (with 20 balls remaining)")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# Row 3 - the original match row, shifted down and right by the new matchNo column
$row3 = @("29th","Punjab Kings","Chris Jordan","c Lalit Yadav b Rabada","2","3","0","0","66.66","Delhi Capitals","Ahmedabad","May 02","Capitals won by 7 wickets (with 14 balls remaining)")
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
